# Remove the obsolete "Change old style hashes to new hashes" TODO item.
# (The remaining hash-related item was reworded/replaced elsewhere in the
# project, so this now-stale bullet point is deleted in its entirety,
# including its paragraph mark, so the list collapses cleanly around it.)

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Change old style hashes to new hashes", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # $rng now covers just the matched text; grab the paragraph that holds it
    # so the delete also removes the trailing paragraph mark (i.e. the whole
    # bullet item disappears instead of leaving an empty one behind).
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
}
